$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Weekly Quantity": append a new weekly-quantity row
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Cells.Item(8, 1).NumberFormat = $wsWeekly.Cells.Item(7, 1).NumberFormat
$wsWeekly.Cells.Item(8, 1).Value = 45662.99999999999
$wsWeekly.Cells.Item(8, 2).Value = 10

# ---------------------------------------------------------------------------
# Sheet "Monthly Trend": append a new monthly-quantity row
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(6, 1).NumberFormat = $wsMonthly.Cells.Item(5, 1).NumberFormat
$wsMonthly.Cells.Item(6, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(6, 2).Value = 10

# ---------------------------------------------------------------------------
# Sheet "PO Forecast": refreshed forecast model - update existing rows and
# append one new forecast row
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$wsForecast.Cells.Item(2, 2).Value = 27
$wsForecast.Cells.Item(3, 2).Value = 27
$wsForecast.Cells.Item(4, 2).Value = 28
$wsForecast.Cells.Item(5, 2).Value = 46
$wsForecast.Cells.Item(6, 2).Value = 53
$wsForecast.Cells.Item(7, 2).Value = 61

$wsForecast.Cells.Item(8, 1).Value = 45662.99999999999
$wsForecast.Cells.Item(8, 2).Value = 63

$wsForecast.Cells.Item(9, 1).Value = 45669.99999999999
$wsForecast.Cells.Item(9, 2).Value = 63

$wsForecast.Cells.Item(10, 1).Value = 45676.99999999999
$wsForecast.Cells.Item(10, 2).Value = 64

$wsForecast.Cells.Item(11, 1).Value = 45683.99999999999
$wsForecast.Cells.Item(11, 2).Value = 64

$wsForecast.Cells.Item(12, 1).Value = 45690.99999999999
$wsForecast.Cells.Item(12, 2).Value = 65

$wsForecast.Cells.Item(13, 1).Value = 45697.99999999999
$wsForecast.Cells.Item(13, 2).Value = 65

$wsForecast.Cells.Item(14, 1).Value = 45704.99999999999
$wsForecast.Cells.Item(14, 2).Value = 66

$wsForecast.Cells.Item(15, 1).Value = 45711.99999999999
$wsForecast.Cells.Item(15, 2).Value = 66

$wsForecast.Cells.Item(16, 1).NumberFormat = $wsForecast.Cells.Item(15, 1).NumberFormat
$wsForecast.Cells.Item(16, 1).Value = 45718.99999999999
$wsForecast.Cells.Item(16, 2).Value = 67
